# Weekly fruit/vegetable price update: a new week's record is inserted as
# row 138 (all existing rows from 138 downward shift down by one), pushing
# the last existing record (old row 176) down to the new row 177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 138, shifting rows 138:176 down to 139:177.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new week's data.
$row = 138
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44463
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112017
$ws.Cells.Item($row, 7).Value = "Apio"
$ws.Cells.Item($row, 8).Value = "Americana (o)"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 40
$ws.Cells.Item($row, 11).Value = 10000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 10000
$ws.Cells.Item($row, 14).Value = "`$/docena de matas"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 1667
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
